$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 123.57143
$ws.Range("J2").Value = 150
$ws.Range("L2").Value = 150
$ws.Range("N2").Value = -376

$ws.Range("H17").Value = 361424.6
$ws.Range("J17").Value = 361424.6
$ws.Range("L17").Value = 1084273.8
$ws.Range("N17").Value = -1084609.8

$ws.Range("H86").Value = 1420.8
$ws.Range("I86").Value = 1166.6666
$ws.Range("J86").Value = 1802
$ws.Range("K86").Value = 1166.6666
$ws.Range("L86").Value = 1802
$ws.Range("M86").Value = -43.66660000000002
$ws.Range("N86").Value = -4048

$ws.Range("H89").Value = 1420.8
$ws.Range("I89").Value = 1166.6666
$ws.Range("J89").Value = 1802
$ws.Range("K89").Value = 5833.333000000001
$ws.Range("L89").Value = 9010
$ws.Range("M89").Value = -217.3330000000005
$ws.Range("N89").Value = -20242

$ws.Range("H101").Value = 66666784
$ws.Range("I101").Value = 100000080
$ws.Range("J101").Value = 185
$ws.Range("K101").Value = 300000240
$ws.Range("L101").Value = 555
$ws.Range("M101").Value = -299998618
$ws.Range("N101").Value = -3799

$ws.Range("H113").Value = 1532.0834
$ws.Range("I113").Value = 1570
$ws.Range("J113").Value = 1505
$ws.Range("K113").Value = 1570
$ws.Range("L113").Value = 1505
$ws.Range("M113").Value = 1684
$ws.Range("N113").Value = -8013

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1182598
$ws.Range("I32").Value = 1395103.1
$ws.Range("K32").Value = 1395103.1
$ws.Range("M32").Value = -1394816.1

$ws.Range("H74").Value = 19657.963
$ws.Range("I74").Value = 28177.082
$ws.Range("J74").Value = 1116.3529
$ws.Range("K74").Value = 28177.082
$ws.Range("L74").Value = 1116.3529
$ws.Range("M74").Value = -27303.082
$ws.Range("N74").Value = -2864.3529

$ws.Range("H77").Value = 19657.963
$ws.Range("I77").Value = 28177.082
$ws.Range("J77").Value = 1116.3529
$ws.Range("K77").Value = 140885.41
$ws.Range("L77").Value = 5581.7645
$ws.Range("M77").Value = -136517.41
$ws.Range("N77").Value = -14317.7645

$ws.Range("H132").Value = 2245271.2
$ws.Range("I132").Value = 2553656.8
$ws.Range("J132").Value = 1011729.2
$ws.Range("K132").Value = 7660970.399999999
$ws.Range("L132").Value = 3035187.6
$ws.Range("M132").Value = -7658440.399999999
$ws.Range("N132").Value = -3040247.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 673.35297
$ws.Range("I94").Value = 333
$ws.Range("J94").Value = 1297.3334
$ws.Range("K94").Value = 333
$ws.Range("L94").Value = 1297.3334
$ws.Range("M94").Value = 118
$ws.Range("N94").Value = -2199.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 2133.3333
$ws.Range("I17").Value = 500
$ws.Range("K17").Value = 500
$ws.Range("M17").Value = -326

$ws.Range("H31").Value = 8134.8965
$ws.Range("I31").Value = 6259.3477
$ws.Range("J31").Value = 15324.5
$ws.Range("K31").Value = 6259.3477
$ws.Range("L31").Value = 15324.5
$ws.Range("M31").Value = -5964.3477
$ws.Range("N31").Value = -15914.5

$ws.Range("H34").Value = 8134.8965
$ws.Range("I34").Value = 6259.3477
$ws.Range("J34").Value = 15324.5
$ws.Range("K34").Value = 6259.3477
$ws.Range("L34").Value = 15324.5
$ws.Range("M34").Value = -6057.3477
$ws.Range("N34").Value = -15728.5

$ws.Range("H99").Value = 4650
$ws.Range("I99").Value = 4657.143
$ws.Range("J99").Value = 4640
$ws.Range("K99").Value = 4657.143
$ws.Range("L99").Value = 4640
$ws.Range("M99").Value = -3159.143
$ws.Range("N99").Value = -7636

$ws.Range("H126").Value = 4650
$ws.Range("I126").Value = 4657.143
$ws.Range("J126").Value = 4640
$ws.Range("K126").Value = 13971.429
$ws.Range("L126").Value = 13920
$ws.Range("M126").Value = -11501.429
$ws.Range("N126").Value = -18860

$ws.Range("H132").Value = 1377.7142
$ws.Range("I132").Value = 918.86365
$ws.Range("J132").Value = 2154.2307
$ws.Range("K132").Value = 2756.59095
$ws.Range("L132").Value = 6462.6921
$ws.Range("M132").Value = -226.5909499999998
$ws.Range("N132").Value = -11522.6921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 29762756
$ws.Range("J131").Value = 34723130
$ws.Range("L131").Value = 104169390
$ws.Range("N131").Value = -104179470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 805.619
$ws.Range("I97").Value = 734.8333
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 734.8333
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -238.8333
$ws.Range("N97").Value = -1892

$ws.Range("H113").Value = 940.9
$ws.Range("I113").Value = 1144.2
$ws.Range("J113").Value = 737.6
$ws.Range("K113").Value = 1144.2
$ws.Range("L113").Value = 737.6
$ws.Range("M113").Value = 1025.8
$ws.Range("N113").Value = -5077.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 28875.514
$ws.Range("I40").Value = 1192.5416
$ws.Range("J40").Value = 79982.53999999999
$ws.Range("K40").Value = 1192.5416
$ws.Range("L40").Value = 79982.53999999999
$ws.Range("M40").Value = -1056.5416
$ws.Range("N40").Value = -80254.53999999999

$ws.Range("H132").Value = 296617.5
$ws.Range("I132").Value = 69225.164
$ws.Range("J132").Value = 916778.4
$ws.Range("K132").Value = 207675.492
$ws.Range("L132").Value = 2750335.2
$ws.Range("M132").Value = -205145.492
$ws.Range("N132").Value = -2755395.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3809.7778
$ws.Range("I62").Value = 3258
$ws.Range("J62").Value = 4499.5
$ws.Range("K62").Value = 3258
$ws.Range("L62").Value = 4499.5
$ws.Range("M62").Value = -2634
$ws.Range("N62").Value = -5747.5

$ws.Range("H65").Value = 3809.7778
$ws.Range("I65").Value = 3258
$ws.Range("J65").Value = 4499.5
$ws.Range("K65").Value = 16290
$ws.Range("L65").Value = 22497.5
$ws.Range("M65").Value = -13170
$ws.Range("N65").Value = -28737.5

$ws.Range("H107").Value = 336.625
$ws.Range("I107").Value = 357.2
$ws.Range("J107").Value = 302.33334
$ws.Range("K107").Value = 1071.6
$ws.Range("L107").Value = 907.0000200000001
$ws.Range("M107").Value = 848.4000000000001
$ws.Range("N107").Value = -4747.00002

$ws.Range("H122").Value = 4387.815
$ws.Range("I122").Value = 1782
$ws.Range("J122").Value = 9599.444
$ws.Range("K122").Value = 5346
$ws.Range("L122").Value = 28798.332
$ws.Range("M122").Value = -2896
$ws.Range("N122").Value = -33698.33199999999

$ws.Range("H132").Value = 3280.2888
$ws.Range("I132").Value = 940.25
$ws.Range("J132").Value = 5152.32
$ws.Range("K132").Value = 2820.75
$ws.Range("L132").Value = 15456.96
$ws.Range("M132").Value = -290.75
$ws.Range("N132").Value = -20516.96

$ws.Range("H136").Value = 1131927.4
$ws.Range("I136").Value = 1348694.9
$ws.Range("K136").Value = 4046084.7
$ws.Range("M136").Value = -4043534.7
